$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 5-11 (Fecha=D, Volumen=M, Precio minimo=N,
# Precio maximo=O, Precio promedio ponderado=P, Precio $/Kg=S)
$data = @(
    @{Row=5;  D=44445; M=160; N=14000; O=15000; P=14500; S=7250},
    @{Row=6;  D=44446; M=300; N=14000; O=15000; P=14500; S=7250},
    @{Row=7;  D=44452; M=200; N=13000; O=14000; P=13500; S=6750},
    @{Row=8;  D=44448; M=100; N=14000; O=15000; P=14500; S=7250},
    @{Row=9;  D=44455; M=160; N=13000; O=14000; P=13500; S=6750},
    @{Row=10; D=44466; M=160; N=13500; O=14000; P=13750; S=6875},
    @{Row=11; D=44468; M=300; N=13000; O=14000; P=13500; S=6750}
)

foreach ($row in $data) {
    $r = $row.Row
    $ws.Cells.Item($r, 4).Value2  = $row.D   # Column D - Fecha
    $ws.Cells.Item($r, 13).Value2 = $row.M   # Column M - Volumen
    $ws.Cells.Item($r, 14).Value2 = $row.N   # Column N - Precio minimo
    $ws.Cells.Item($r, 15).Value2 = $row.O   # Column O - Precio maximo
    $ws.Cells.Item($r, 16).Value2 = $row.P   # Column P - Precio promedio ponderado
    $ws.Cells.Item($r, 19).Value2 = $row.S   # Column S - Precio $/Kg
}
